$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Relocate the "_GoBack" bookmark.
#    In the original document it wraps the final screenshot picture
#    (the last thing that was edited in a previous session). The new
#    revision instead has the user's last edit point sitting inside the
#    "(day2)" heading text, between "day2" and the closing paren - so the
#    bookmark needs to move there. Word always keeps a single "_GoBack"
#    bookmark, so remove the old one first, then add the new one; adding a
#    bookmark on a collapsed (zero-length) range in the middle of a run
#    causes Word to split that run in two around the bookmark, which is
#    exactly the structure the target document has.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$content = $d.Content.Text
$day2Idx = $content.IndexOf("(day2)")
if ($day2Idx -ge 0) {
    $splitAt = $day2Idx + 5   # right after "(day2", before the closing ")"
    $goBackRange = $d.Range($splitAt, $splitAt)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}

# ---------------------------------------------------------------------------
# 2. Split the run that holds the literal text ".tsv (the closing smart
#    quote + period + the word "tsv") into two runs - one ending right
#    after the period, and one holding just "tsv" - mirroring the word
#    boundary Word's proofing pass breaks the text on. Re-fetch the text
#    each time since the previous step may have shifted offsets.
# ---------------------------------------------------------------------------
$content = $d.Content.Text
$tsvIdx = $content.IndexOf([char]0x201C + ".tsv")
if ($tsvIdx -lt 0) {
    $tsvIdx = $content.IndexOf("“.tsv")
}
if ($tsvIdx -ge 0) {
    $splitAt = $tsvIdx + 2   # after the opening smart-quote + ".", before "tsv"
    $tempRange = $d.Range($splitAt, $splitAt)
    $d.Bookmarks.Add("ZZTempSplit", $tempRange)
    $d.Bookmarks("ZZTempSplit").Delete()
}

Write-Output "done"
